$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen the data range and write the new access-log rows (101-132).
# Apply a temporary text number format to the target range so that
# date-like ("2025-06-04") and time-like ("15:28:20") strings are stored
# as literal text instead of being auto-converted to date/time serials,
# matching the inlineStr cells used by the rest of the sheet. The format
# is cleared again afterwards so the cells end up with no explicit style,
# just like the appended rows in the source diff.
$newRange = $ws.Range("A101:C132")
$newRange.NumberFormat = "@"

$ws.Range("A101").Value = "coord123"
$ws.Range("B101").Value = "2025-06-04"
$ws.Range("C101").Value = "15:28:20"
$ws.Range("A102").Value = "coord123"
$ws.Range("B102").Value = "2025-06-04"
$ws.Range("C102").Value = "15:32:02"
$ws.Range("A103").Value = "coord123"
$ws.Range("B103").Value = "2025-06-04"
$ws.Range("C103").Value = "15:33:20"
$ws.Range("A104").Value = "coord123"
$ws.Range("B104").Value = "2025-06-04"
$ws.Range("C104").Value = "15:34:53"
$ws.Range("A105").Value = "profana"
$ws.Range("B105").Value = "2025-06-04"
$ws.Range("C105").Value = "15:37:15"
$ws.Range("A106").Value = "coord123"
$ws.Range("B106").Value = "2025-06-04"
$ws.Range("C106").Value = "15:37:50"
$ws.Range("A107").Value = "coord123"
$ws.Range("B107").Value = "2025-06-04"
$ws.Range("C107").Value = "15:50:52"
$ws.Range("A108").Value = "coord123"
$ws.Range("B108").Value = "2025-06-04"
$ws.Range("C108").Value = "16:01:40"
$ws.Range("A109").Value = "coord123"
$ws.Range("B109").Value = "2025-06-04"
$ws.Range("C109").Value = "16:04:00"
$ws.Range("A110").Value = "coord123"
$ws.Range("B110").Value = "2025-06-04"
$ws.Range("C110").Value = "16:07:05"
$ws.Range("A111").Value = "coord123"
$ws.Range("B111").Value = "2025-06-04"
$ws.Range("C111").Value = "16:07:41"
$ws.Range("A112").Value = "profana"
$ws.Range("B112").Value = "2025-06-04"
$ws.Range("C112").Value = "16:14:56"
$ws.Range("A113").Value = "profana"
$ws.Range("B113").Value = "2025-06-04"
$ws.Range("C113").Value = "16:16:27"
$ws.Range("A114").Value = "profana"
$ws.Range("B114").Value = "2025-06-04"
$ws.Range("C114").Value = "16:17:05"
$ws.Range("A115").Value = "profana"
$ws.Range("B115").Value = "2025-06-04"
$ws.Range("C115").Value = "16:17:55"
$ws.Range("A116").Value = "profana"
$ws.Range("B116").Value = "2025-06-04"
$ws.Range("C116").Value = "16:18:03"
$ws.Range("A117").Value = "profana"
$ws.Range("B117").Value = "2025-06-04"
$ws.Range("C117").Value = "16:18:48"
$ws.Range("A118").Value = "profana"
$ws.Range("B118").Value = "2025-06-04"
$ws.Range("C118").Value = "16:31:35"
$ws.Range("A119").Value = "profana"
$ws.Range("B119").Value = "2025-06-04"
$ws.Range("C119").Value = "16:33:09"
$ws.Range("A120").Value = "profana"
$ws.Range("B120").Value = "2025-06-04"
$ws.Range("C120").Value = "16:33:37"
$ws.Range("A121").Value = "profana"
$ws.Range("B121").Value = "2025-06-04"
$ws.Range("C121").Value = "16:41:41"
$ws.Range("A122").Value = "profana"
$ws.Range("B122").Value = "2025-06-04"
$ws.Range("C122").Value = "16:44:30"
$ws.Range("A123").Value = "profana"
$ws.Range("B123").Value = "2025-06-04"
$ws.Range("C123").Value = "16:46:57"
$ws.Range("A124").Value = "profana"
$ws.Range("B124").Value = "2025-06-04"
$ws.Range("C124").Value = "16:51:07"
$ws.Range("A125").Value = "coord123"
$ws.Range("B125").Value = "2025-06-04"
$ws.Range("C125").Value = "16:52:39"
$ws.Range("A126").Value = "coord123"
$ws.Range("B126").Value = "2025-06-04"
$ws.Range("C126").Value = "16:53:36"
$ws.Range("A127").Value = "coord123"
$ws.Range("B127").Value = "2025-06-04"
$ws.Range("C127").Value = "16:54:32"
$ws.Range("A128").Value = "coord123"
$ws.Range("B128").Value = "2025-06-04"
$ws.Range("C128").Value = "16:54:55"
$ws.Range("A129").Value = "coord123"
$ws.Range("B129").Value = "2025-06-04"
$ws.Range("C129").Value = "18:45:28"
$ws.Range("A130").Value = "coord123"
$ws.Range("B130").Value = "2025-06-04"
$ws.Range("C130").Value = "18:46:05"
$ws.Range("A131").Value = "profana"
$ws.Range("B131").Value = "2025-06-04"
$ws.Range("C131").Value = "18:51:30"
$ws.Range("A132").Value = "coord123"
$ws.Range("B132").Value = "2025-06-04"
$ws.Range("C132").Value = "18:58:37"

$newRange.ClearFormats()

